$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOLX")

# Row 4 - Inventory
$ws.Range("B4").Value = 456000000.0
$ws.Range("C4").Value = 420000000.0
$ws.Range("D4").Value = 395000000.0
$ws.Range("E4").Value = 414000000.0
$ws.Range("F4").Value = 401000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 225000000.0
$ws.Range("C14").Value = 184000000.0
$ws.Range("D14").Value = 179000000.0
$ws.Range("E14").Value = 127000000.0
$ws.Range("F14").Value = 134000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("B24").Value = 230000000.0
$ws.Range("C24").Value = 189000000.0
$ws.Range("D24").Value = 186000000.0
$ws.Range("E24").Value = 214000000.0
$ws.Range("F24").Value = 236000000.0

# Row 37 - Net Debt
$ws.Range("G37").Value = 2719600000.0

# Row 38 - Total Debt
$ws.Range("G38").Value = 3090400000.0

$wb.Save()
